# Updates the cryptos price/volume table to the values from the latest
# GitHub Actions scrape. Price cells that look like plain numbers (e.g.
# "574.60", "0.0534") are written with a leading apostrophe so Excel keeps
# them as text (matching the workbook's inline-string cells) instead of
# silently converting them to numeric values and dropping trailing zeros.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '62.918.07'
$ws.Range('E2').Value = '  -1.44%  '
$ws.Range('D3').Value = '2.544.91'
$ws.Range('E3').Value = '  +0.00%  '
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').Value = '''574.60'
$ws.Range('E5').Value = '  +0.03%  '
$ws.Range('D6').Value = '''146.07'
$ws.Range('E6').Value = '  -1.61%  '
$ws.Range('E7').Value = '  +0.00%  '
$ws.Range('E8').Value = '  -1.44%  '
$ws.Range('E9').Value = '  -1.82%  '
$ws.Range('E10').Value = '  -4.73%  '
$ws.Range('E11').Value = '  -0.55%  '
$ws.Range('D12').Value = '''0.353'
$ws.Range('E12').Value = '  -1.28%  '
$ws.Range('D13').Value = '''27.06'
$ws.Range('E13').Value = '  -4.11%  '
$ws.Range('D14').Value = '2.999.01'
$ws.Range('E14').Value = '  -0.05%  '
$ws.Range('D15').Value = '62.841.95'
$ws.Range('E15').Value = '  -1.21%  '
$ws.Range('E16').Value = '  -1.46%  '
$ws.Range('D17').Value = '2.546.64'
$ws.Range('E17').Value = '  +0.09%  '
$ws.Range('E18').Value = '  -2.15%  '
$ws.Range('D19').Value = '''335.12'
$ws.Range('E19').Value = '  -1.91%  '
$ws.Range('E20').Value = '  -1.32%  '
$ws.Range('D21').Value = '''6.74'
$ws.Range('E21').Value = '  -2.22%  '
$ws.Range('E22').Value = '  +0.13%  '
$ws.Range('D23').Value = '''65.16'
$ws.Range('E23').Value = '  -1.64%  '
$ws.Range('E24').Value = '  -0.62%  '
$ws.Range('E25').Value = '  +1.48%  '
$ws.Range('E26').Value = '  +0.06%  '
$ws.Range('D27').Value = '''8.31'
$ws.Range('E27').Value = '  -0.66%  '
$ws.Range('D28').Value = '''1.46'
$ws.Range('E28').Value = '  +2.95%  '
$ws.Range('D29').Value = '''7.24'
$ws.Range('E29').Value = '  +4.67%  '
$ws.Range('E30').Value = '  -3.27%  '
$ws.Range('D31').Value = '''1.86'
$ws.Range('E31').Value = '  -0.66%  '
$ws.Range('D32').Value = '''177.71'
$ws.Range('E32').Value = '  -0.21%  '
$ws.Range('D33').Value = '''1.54'
$ws.Range('E33').Value = '  -3.80%  '
$ws.Range('D34').Value = '''404.04'
$ws.Range('E34').Value = '  -4.20%  '
$ws.Range('D35').Value = '''19.07'
$ws.Range('E35').Value = '  -0.25%  '
$ws.Range('E36').Value = '  -1.99%  '
$ws.Range('D38').Value = '''4.32'
$ws.Range('E38').Value = '  -2.22%  '
$ws.Range('D39').Value = '''1.72'
$ws.Range('E39').Value = '  -2.09%  '
$ws.Range('D40').Value = '''0.999'
$ws.Range('E40').Value = '  -0.12%  '
$ws.Range('D41').Value = '''39.32'
$ws.Range('E41').Value = '  -3.45%  '
$ws.Range('D42').Value = '''150.82'
$ws.Range('E42').Value = '  -1.68%  '
$ws.Range('E43').Value = '  -1.68%  '
$ws.Range('D44').Value = '''20.72'
$ws.Range('D45').Value = '''0.0534'
$ws.Range('E45').Value = '  +0.09%  '
$ws.Range('D46').Value = '''0.599'
$ws.Range('E46').Value = '  -2.12%  '
$ws.Range('E47').Value = '  -0.85%  '
$ws.Range('D48').Value = '''0.0238'
$ws.Range('E48').Value = '  +2.29%  '
$ws.Range('D49').Value = '''18.11'
$ws.Range('E49').Value = '  -3.56%  '
$ws.Range('D50').Value = '''11.30'
$ws.Range('E50').Value = '  +0.43%  '
$ws.Range('E51').Value = '  -8.61%  '
